# Add two new columns, I ("I0") and J ("IF"), to the sheet.
# Header cells (row 1) get the same formatting/style as the existing
# header cells (e.g. H1) -- copy format only (xlPasteFormats = -4122) so
# the new cells reuse the existing style index rather than minting a new
# (near-identical) one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# -- Header row (row 1) --------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial($xlPasteFormats)

$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# -- Data rows 2-11 -------------------------------------------------------
$iValues = @{2 = 6; 3 = 7; 4 = 9; 5 = 6; 6 = 7; 7 = 1; 8 = 1; 9 = 1; 10 = 1; 11 = 1}
$jValues = @{2 = 9; 3 = 9; 4 = 9; 5 = 8; 6 = 9; 7 = 7; 8 = 5; 9 = 6; 10 = 4; 11 = 5}

foreach ($row in 2..11) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
